# Apply the "Trade #16 closed" update to the live trading results workbook.
#
# Summary:
#  - Sheet "Summary": update Current Capital, Total P&L $, Total P&L %,
#    Total Trades, Losing Trades, Win Rate %
#  - Sheet "Strategy Status": update the MarketMaking row (Capital, Trades,
#    P&L $, P&L %, Win Rate %)
#  - Sheet "All Trades" and "MarketMaking": append a new trade record (row 17)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02   # Current Capital
$summary.Range("B4").Value = 0.02      # Total P&L $
$summary.Range("B5").Value = 0.02      # Total P&L %
$summary.Range("B6").Value = 16        # Total Trades
$summary.Range("B8").Value = 6         # Losing Trades
$summary.Range("B9").Value = 31.25     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking is row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 16         # Trades
$status.Range("E4").Value = 0.02       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 31.25      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade (#16, array index row 17) to both the
#    "All Trades" sheet and the "MarketMaking" sheet.
# ---------------------------------------------------------------------
$newTradeRow = @{
    A = 16
    B = "2026-02-17"
    C = "04:07:40"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.82
    G = 0.8100000000000001
    H = "CLOSED"
    I = -1.2195
    J = -0.01
    K = 100.02
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.12
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(17, 1).Value = $newTradeRow.A

    # Columns B (date) and C (time) look like dates/times, so Excel would
    # normally auto-convert them to date/time serial numbers on entry. The
    # source data stores them as plain text, so force a Text number format
    # before assigning, then clear the format again so no residual style
    # is left behind on the cell (matching the rest of the sheet).
    $ws.Cells.Item(17, 2).NumberFormat = "@"
    $ws.Cells.Item(17, 2).Value = $newTradeRow.B
    $ws.Cells.Item(17, 2).ClearFormats()

    $ws.Cells.Item(17, 3).NumberFormat = "@"
    $ws.Cells.Item(17, 3).Value = $newTradeRow.C
    $ws.Cells.Item(17, 3).ClearFormats()

    $ws.Cells.Item(17, 4).Value = $newTradeRow.D
    $ws.Cells.Item(17, 5).Value = $newTradeRow.E
    $ws.Cells.Item(17, 6).Value = $newTradeRow.F
    $ws.Cells.Item(17, 7).Value = $newTradeRow.G
    $ws.Cells.Item(17, 8).Value = $newTradeRow.H
    $ws.Cells.Item(17, 9).Value = $newTradeRow.I
    $ws.Cells.Item(17, 10).Value = $newTradeRow.J
    $ws.Cells.Item(17, 11).Value = $newTradeRow.K
    $ws.Cells.Item(17, 12).Value = $newTradeRow.L
    $ws.Cells.Item(17, 13).Value = $newTradeRow.M
    $ws.Cells.Item(17, 14).Value = $newTradeRow.N
    $ws.Cells.Item(17, 15).Value = $newTradeRow.O
    $ws.Cells.Item(17, 16).Value = $newTradeRow.P
    $ws.Cells.Item(17, 17).Value = $newTradeRow.Q
}
